$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data
$ws.Range("A5").Value = "BinaryTree"
$ws.Range("B5").Value = "Ongoing"

# Resize the table (ListObject) to include the new row
$table = $ws.ListObjects.Item(1)
$table.Resize($ws.Range("A1:B5"))

# Extend the color-scale conditional formatting over column B to include the new row
$fc = $ws.Range("B2").FormatConditions
$cond = $fc.Item(2)
$cond.ModifyAppliesToRange($ws.Range("B2:B5"))

# Update selection to match target state
$ws.Range("D8").Select()
